# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
#
# All data cells on the sheet are plain text ("inline string") cells - prices such
# as "301.05" or "1.00" are display strings, not numbers. Assigning a numeric-looking
# string straight to Range.Value lets Excel's type inference silently convert it to a
# real number (losing trailing zeros / thousands-dot formatting), so for column D
# (Price) we momentarily force a Text number format, write the string, then restore
# the default "Normal" style so the cell ends up styled exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Column D (Price) updates ---
Set-TextValue "D2" '43.115.91'
Set-TextValue "D3" '2.305.73'
Set-TextValue "D5" '301.05'
Set-TextValue "D6" '97.69'
Set-TextValue "D9" '0.517'
Set-TextValue "D10" '35.83'
Set-TextValue "D13" '17.93'
Set-TextValue "D14" '6.88'
Set-TextValue "D15" '2.664.89'
Set-TextValue "D16" '2.279.05'
Set-TextValue "D18" '43.002.51'
Set-TextValue "D19" '13.06'
Set-TextValue "D22" '68.32'
Set-TextValue "D23" '237.99'
Set-TextValue "D24" '2.23'
Set-TextValue "D25" '0.990'
Set-TextValue "D28" '25.22'
Set-TextValue "D29" '166.72'
Set-TextValue "D30" '9.15'
Set-TextValue "D32" '33.08'
Set-TextValue "D33" '1.00'
Set-TextValue "D34" '5.14'
Set-TextValue "D35" '18.22'
Set-TextValue "D38" '0.0690'
Set-TextValue "D40" '1.78'
Set-TextValue "D43" '2.008.90'
Set-TextValue "D44" '0.0288'
Set-TextValue "D45" '2.18'
Set-TextValue "D47" '17.49'
Set-TextValue "D49" '54.45'
Set-TextValue "D50" '2.536.92'

# --- Rows 34/35: Filecoin and Celestia swapped ranking positions ---
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("E16").Value = '  -2.21%  '
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("E31").Value = '  -13.36%  '
$ws.Range("E32").Value = '  -5.06%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("E35").Value = '  +1.98%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("E42").Value = '  -2.83%  '
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("E45").Value = '  -7.21%  '
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  -1.24%  '
